$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11").Value = "Max"
$ws.Range("C11").Value = "max@x.com"
$ws.Hyperlinks.Add($ws.Range("C11"), "mailto:max@x.com")
$ws.Range("C11").Style = "Hyperlink"

$ws.Range("B12").Value = "Kat"
$ws.Range("C12").Value = "kat@x.com"
$ws.Hyperlinks.Add($ws.Range("C12"), "mailto:kat@x.com")
$ws.Range("C12").Style = "Hyperlink"

$ws.Range("B13").Select()
